# Apply the workbook edit described by the diff:
#  - Rename Sheet1 -> Home
#  - Add a new worksheet "Manage Units" right after Home, populated with a
#    new scenario/expected-result table, and make it the active tab.
#  - Shrink the row-5 height on Home (content there got shorter) and update
#    the remembered selections on both sheets.

$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet ---
$home = $wb.Worksheets.Item(1)
$home.Name = "Home"

# --- Add the new sheet right after Home (becomes the active sheet) ---
$manage = $wb.Worksheets.Add($null, $home)
$manage.Name = "Manage Units"

# --- Populate "Manage Units" sheet ---
$manage.Range("A1").Value = "Scenario"
$manage.Range("B1").Value = "Expected result"

$manage.Range("A2").Value = "When a user adds a new unit"
$manage.Range("B2").Value = "The new unit should be displayed in the Units table"

$manage.Range("A3").Value = "When a user attempts to adds more than 7 alphanumeric characters to the unitcode"
$manage.Range("B3").Value = "The text should not exceed 7 characters"

$manage.Range("A4").Value = "When the user enters less than 7 characters in the unitcode"
$manage.Range("B4").Value = "The user should be prompted to add 7 characters"

$manage.Range("A5").Value = "When the user adds a unit title to a unit"
$manage.Range("B5").Value = "The units table should reflect the unit with that unit title"

$manage.Range("A6").Value = "When the user "

# --- Row heights on Manage Units sheet ---
$manage.Rows.Item(1).RowHeight = 28.5
$manage.Rows.Item(2).RowHeight = 99.75
$manage.Rows.Item(3).RowHeight = 28.5
$manage.Rows.Item(4).RowHeight = 28.5
$manage.Rows.Item(5).RowHeight = 28.5

# --- Update existing Home sheet: row 5 height shrinks (text shortened) ---
$home.Rows.Item(5).RowHeight = 28.5

# --- Selections remembered per-sheet ---
$home.Range("A4").Select()
$manage.Range("C4").Select()

# --- "Manage Units" ends up as the active/selected tab ---
$manage.Activate()
$manage.Range("C4").Select()
